# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 372
$ws1.Range("F3").Value = 800
$ws1.Range("F5").Value = 905
$ws1.Range("F6").Value = 2198
$ws1.Range("F7").Value = 195

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 372
$ws4.Range("F3").Value = 800
$ws4.Range("F7").Value = 905
$ws4.Range("F8").Value = 2198
$ws4.Range("F10").Value = 195
